$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (player, position(s), team) replacing/extending the existing
# A2:C18 range and adding a new row 19.
$data = @(
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Nick Richards", "C", "Phoenix Suns"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Chicago Bulls"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
